# Apply weekly fruit/vegetable price update: reorder the per-date values
# (date, volume, min/max/avg price, price per kg) across rows 2-9.
# Row 5 is unchanged; the other rows' D/M/N/O/P/S values are shuffled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = 44193

$ws.Range("D3").Value2 = 44188
$ws.Range("N3").Value2 = 15000
$ws.Range("O3").Value2 = 15000
$ws.Range("P3").Value2 = 15000
$ws.Range("S3").Value2 = 3000

$ws.Range("D4").Value2 = 44189
$ws.Range("M4").Value2 = 40
$ws.Range("N4").Value2 = 15000
$ws.Range("O4").Value2 = 15000
$ws.Range("P4").Value2 = 15000
$ws.Range("S4").Value2 = 3000

$ws.Range("D6").Value2 = 44186
$ws.Range("M6").Value2 = 40
$ws.Range("N6").Value2 = 15000
$ws.Range("O6").Value2 = 15000
$ws.Range("P6").Value2 = 15000
$ws.Range("S6").Value2 = 3000

$ws.Range("D7").Value2 = 44175
$ws.Range("M7").Value2 = 25
$ws.Range("N7").Value2 = 20000
$ws.Range("O7").Value2 = 20000
$ws.Range("P7").Value2 = 20000
$ws.Range("S7").Value2 = 4000

$ws.Range("D8").Value2 = 44181
$ws.Range("M8").Value2 = 30
$ws.Range("N8").Value2 = 20000
$ws.Range("O8").Value2 = 20000
$ws.Range("P8").Value2 = 20000
$ws.Range("S8").Value2 = 4000

$ws.Range("D9").Value2 = 44179
$ws.Range("M9").Value2 = 45
$ws.Range("N9").Value2 = 20000
$ws.Range("O9").Value2 = 20000
$ws.Range("P9").Value2 = 20000
$ws.Range("S9").Value2 = 4000
